$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2 through 45 is updated from
# serial date 45204 (2023-10-05) to 45205 (2023-10-06).
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
